$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Update the "twelve months ended" column headers: the reporting window
# rolls forward by one fiscal year (drop 1396/12, keep 1397-1400/12, add 1401/12) ---
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Update the data rows (values also shift one year left, with the newest
# year, 1401/12, populated with new figures) ---

# row 10: هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 5110
$ws.Range("F10").Value = 2628
$ws.Range("G10").Value = 3478
$ws.Range("H10").Value = 3967
$ws.Range("I10").Value = 6390

# row 13: هزینه تبلیغات
$ws.Range("E13").Value = 23690
$ws.Range("F13").Value = 12673
$ws.Range("G13").Value = 7238
$ws.Range("H13").Value = 63602
$ws.Range("I13").Value = 58900

# row 14: هزینه مواد مصرفی
$ws.Range("E14").Value = 3233
$ws.Range("F14").Value = 4910
$ws.Range("G14").Value = 6887
$ws.Range("H14").Value = 15801
$ws.Range("I14").Value = 34921

# row 15: هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 95
$ws.Range("F15").Value = 163
$ws.Range("G15").Value = 1038
$ws.Range("H15").Value = 816
$ws.Range("I15").Value = 682

# row 16: هزینه استهلاک
$ws.Range("E16").Value = 7229
$ws.Range("F16").Value = 7585
$ws.Range("G16").Value = 6757
$ws.Range("H16").Value = 6141
$ws.Range("I16").Value = 8631

# row 17: هزینه حقوق و دستمزد
$ws.Range("E17").Value = 81500
$ws.Range("F17").Value = 89264
$ws.Range("G17").Value = 110314
$ws.Range("H17").Value = 148241
$ws.Range("I17").Value = 256312

# row 18: هزینه مطالبات مشکوک الوصول
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# row 19: سایر هزینه ها
$ws.Range("E19").Value = 148889
$ws.Range("F19").Value = 59374
$ws.Range("G19").Value = 57760
$ws.Range("H19").Value = 198919
$ws.Range("I19").Value = 217376

# row 20: جمع (total)
$ws.Range("E20").Value = 269746
$ws.Range("F20").Value = 176597
$ws.Range("G20").Value = 193472
$ws.Range("H20").Value = 437487
$ws.Range("I20").Value = 583212

# row 26: تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 277
$ws.Range("F26").Value = 287
$ws.Range("G26").Value = 187
$ws.Range("H26").Value = 234
$ws.Range("I26").Value = 241

# row 27: تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 129
$ws.Range("F27").Value = 119
$ws.Range("G27").Value = 197
$ws.Range("H27").Value = 147
$ws.Range("I27").Value = 143
